# Fix typo in heading that broke XLS file chunking:
# The header cell G3 ("Description") should be lowercase ("description").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Correct the typo in the column header
$ws.Range("G3").Value = "description"

# Update the active selection to match the saved workbook state
$ws.Range("G4").Select()
